$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D updates: use a Text number format for values that would
# otherwise be auto-coerced to numbers (e.g. "602.00" -> 602), then
# clear the format again so the cell keeps its original (unstyled) look.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.456.75"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.641.54"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "602.00"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "146.34"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "5.60"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "27.33"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "3.122.03"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "63.321.32"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "2.655.59"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "340.98"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "6.92"
$ws.Range("E21").Value = "  +3.22%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "5.56"
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("D24").Value = "66.75"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("D25").Value = "1.69"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "8.88"
$ws.Range("E26").Value = "  +4.79%  "
$ws.Range("D27").Value = "1.53"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "547.53"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").Value = "7.84"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  +6.95%  "
$ws.Range("D36").Value = "167.90"
$ws.Range("E36").Value = "  -4.13%  "
$ws.Range("D37").Value = "0.406"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "19.09"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  +7.07%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "169.31"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "22.52"
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("D45").Value = "0.0578"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("D51").Value = "11.28"
$ws.Range("E51").Value = "  -0.51%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D51").ClearFormats()
